$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C (old C..Q shift right to D..R)
$ws.Columns.Item(3).Insert()

# New column header
$ws.Range("C1").Value = "Term Type"

# New column values for selected data rows
$ws.Range("C4").Value = "germplasm passport"
$ws.Range("C5").Value = "PHENOTYPE"

# Update selection to match the target view state
$ws.Range("C11").Select()
